# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E22) is renumbered from the old
# descending block (1903,1902,1901,1812,1811,1810,1809) to the new
# ascending block (1809,1810,1811,1812,1901,1902,1903), and the "Valor
# Mora" amounts for the first and last periods (F16 / F22) swap so the
# value that used to belong to period 1903 now travels with it to row 22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "1809"
$ws.Range("E17").Value = "1810"
$ws.Range("E18").Value = "1811"
$ws.Range("E19").Value = "1812"
$ws.Range("E20").Value = "1901"
$ws.Range("E21").Value = "1902"
$ws.Range("E22").Value = "1903"

$ws.Range("F16").Value = 31249
$ws.Range("F22").Value = 26041
